$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.228.74"
$ws.Range("E2").Value = "  +3.33%  "
$ws.Range("D3").Value = "2.336.50"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'545.12"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "2.332.68"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "2.748.83"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "60.181.15"
$ws.Range("E16").Value = "  +3.46%  "
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "2.332.03"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").Value = "'10.61"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("D21").Value = "'6.77"
$ws.Range("E21").Value = "  +5.55%  "
$ws.Range("D22").Value = "'313.69"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").Value = "'63.61"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("D25").Value = "'0.170"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'7.90"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("E28").Value = "  +6.74%  "
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("E31").Value = "  +11.56%  "
$ws.Range("D32").Value = "0.0₃0730"
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("E34").Value = "  +12.10%  "
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").Value = "'17.98"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("E39").Value = "  +5.76%  "
$ws.Range("D40").Value = "'321.33"
$ws.Range("E40").Value = "  +11.04%  "
$ws.Range("D41").Value = "'38.13"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("E42").Value = "  +1.87%  "
$ws.Range("D43").Value = "'140.43"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").Value = "'3.46"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D45").Value = "'0.0945"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").Value = "'19.45"
$ws.Range("E46").Value = "  +7.28%  "
$ws.Range("D47").Value = "'0.0497"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("E49").Value = "  +1.01%  "

# Rows 50 and 51 swap coins/links with freshly updated price and volume data
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0213"
$ws.Range("E50").Value = "  +13.51%  "

$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'11.02"
$ws.Range("E51").Value = "  +0.81%  "
